# Fruta / hortaliza, semanal
# Insert two new rows after row 64 so that the current row 65 (Agrícola del
# Norte S.A. de Arica - Guayaba, 2021-10-04) shifts down to row 67, and
# populate the two newly-inserted rows (65 & 66) with the values that used
# to live in rows 63 & 64 before this week's price update. Rows 63 & 64
# themselves are then refreshed with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two blank rows at position 65 (pushes old row 65 -> row 67) ---
$ws.Rows.Item(65).Insert()
$ws.Rows.Item(65).Insert()

# --- New row 65 (previous content of row 63, Calidad = Primera) ---
$ws.Range("A65").Value = 1
$ws.Range("B65").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C65").Value = "Arica y Parinacota"
$ws.Range("D65").Value = 44722
$ws.Range("E65").Value = 15
$ws.Range("F65").Value = "Fruta"
$ws.Range("G65").Value = 100108
$ws.Range("H65").Value = "Tropicales y subtropicales"
$ws.Range("I65").Value = 100108001
$ws.Range("J65").Value = "Guayaba"
$ws.Range("K65").Value = "Sin especificar"
$ws.Range("L65").Value = "Primera"
$ws.Range("M65").Value = 140
$ws.Range("N65").Value = 800
$ws.Range("O65").Value = 900
$ws.Range("P65").Value = 850
$ws.Range("Q65").Value = "$/kilo (en caja de 10 kilos )"
$ws.Range("R65").Value = "Región de Arica y Parinacota"
$ws.Range("S65").Value = 850
$ws.Range("T65").Value = 1

# --- New row 66 (previous content of row 64, Calidad = Segunda) ---
$ws.Range("A66").Value = 1
$ws.Range("B66").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C66").Value = "Arica y Parinacota"
$ws.Range("D66").Value = 44722
$ws.Range("E66").Value = 15
$ws.Range("F66").Value = "Fruta"
$ws.Range("G66").Value = 100108
$ws.Range("H66").Value = "Tropicales y subtropicales"
$ws.Range("I66").Value = 100108001
$ws.Range("J66").Value = "Guayaba"
$ws.Range("K66").Value = "Sin especificar"
$ws.Range("L66").Value = "Segunda"
$ws.Range("M66").Value = 200
$ws.Range("N66").Value = 700
$ws.Range("O66").Value = 800
$ws.Range("P66").Value = 750
$ws.Range("Q66").Value = "$/kilo (en caja de 10 kilos )"
$ws.Range("R66").Value = "Región de Arica y Parinacota"
$ws.Range("S66").Value = 750
$ws.Range("T66").Value = 1

# --- Update row 63 with this week's new values ---
$ws.Range("D63").Value = 45075
$ws.Range("N63").Value = 4500
$ws.Range("O63").Value = 5000
$ws.Range("P63").Value = 4786
$ws.Range("Q63").Value = "$/caja 10 kilos"
$ws.Range("S63").Value = 479
$ws.Range("T63").Value = 10

# --- Update row 64 with this week's new values ---
$ws.Range("D64").Value = 45075
$ws.Range("M64").Value = 140
$ws.Range("N64").Value = 3500
$ws.Range("O64").Value = 4000
$ws.Range("P64").Value = 3643
$ws.Range("Q64").Value = "$/caja 10 kilos"
$ws.Range("S64").Value = 364
$ws.Range("T64").Value = 10
